$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.069782333333333
$ws.Range("H2").Value = 3.209347
$ws.Range("I2").Value = 0.5403049320348338
$ws.Range("J2").Value = 0.5403049320348337
$ws.Range("M2").Value = 52.63198466666667
$ws.Range("N2").Value = 157.895954
$ws.Range("O2").Value = 0.1037886003335349
$ws.Range("P2").Value = 0.1037886003335349
$ws.Range("Q2").Value = 56.3047673646709
$ws.Range("R2").Value = 506.7429062820381
$ws.Range("S2").Value = 0.0560774926492011
$ws.Range("T2").Value = 0.0560774926492011
$ws.Range("A3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.069782333333333
$ws.Range("H3").Value = 3.209347
$ws.Range("I3").Value = 0.5403049320348338
$ws.Range("J3").Value = 0.5403049320348337
$ws.Range("O3").Value = 0.2720850929153589
$ws.Range("P3").Value = 0.2720850929153589
$ws.Range("Q3").Value = 147.6047254781625
$ws.Range("R3").Value = 1328.442529303463
$ws.Range("S3").Value = 0.1470089176353244
$ws.Range("T3").Value = 0.1470089176353244
$ws.Range("A4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.069782333333333
$ws.Range("H4").Value = 3.209347
$ws.Range("I4").Value = 0.5403049320348338
$ws.Range("J4").Value = 0.5403049320348337
$ws.Range("M4").Value = 45.81758366666667
$ws.Range("N4").Value = 137.452751
$ws.Range("O4").Value = 0.09035081822479055
$ws.Range("P4").Value = 0.09035081822479055
$ws.Range("Q4").Value = 49.0148415626219
$ws.Range("R4").Value = 441.1335740635971
$ws.Range("S4").Value = 0.04881699270023707
$ws.Range("T4").Value = 0.04881699270023707
$ws.Range("A5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.069782333333333
$ws.Range("H5").Value = 3.209347
$ws.Range("I5").Value = 0.5403049320348338
$ws.Range("J5").Value = 0.5403049320348337
$ws.Range("M5").Value = 270.6815896666666
$ws.Range("N5").Value = 812.044769
$ws.Range("O5").Value = 0.5337754885263156
$ws.Range("P5").Value = 0.5337754885263157
$ws.Range("Q5").Value = 289.5703825839826
$ws.Range("R5").Value = 2606.133443255843
$ws.Range("S5").Value = 0.2884015290500712
$ws.Range("T5").Value = 0.2884015290500712
$ws.Range("A6").Value = "MuSCs"
$ws.Range("G6").Value = 0.663689
$ws.Range("H6").Value = 1.991067
$ws.Range("I6").Value = 0.3352031799963669
$ws.Range("J6").Value = 0.3352031799963669
$ws.Range("M6").Value = 52.63198466666667
$ws.Range("N6").Value = 157.895954
$ws.Range("O6").Value = 0.1037886003335349
$ws.Range("P6").Value = 0.1037886003335349
$ws.Range("Q6").Value = 34.93126927143534
$ws.Range("R6").Value = 314.381423442918
$ws.Range("S6").Value = 0.03479026887917289
$ws.Range("T6").Value = 0.03479026887917289
$ws.Range("A7").Value = "MuSCs"
$ws.Range("G7").Value = 0.663689
$ws.Range("H7").Value = 1.991067
$ws.Range("I7").Value = 0.3352031799963669
$ws.Range("J7").Value = 0.3352031799963669
$ws.Range("O7").Value = 0.2720850929153589
$ws.Range("P7").Value = 0.2720850929153589
$ws.Range("Q7").Value = 91.57342535526033
$ws.Range("R7").Value = 824.1608281973429
$ws.Range("S7").Value = 0.09120378837483525
$ws.Range("T7").Value = 0.09120378837483527
$ws.Range("A8").Value = "MuSCs"
$ws.Range("G8").Value = 0.663689
$ws.Range("H8").Value = 1.991067
$ws.Range("I8").Value = 0.3352031799963669
$ws.Range("J8").Value = 0.3352031799963669
$ws.Range("M8").Value = 45.81758366666667
$ws.Range("N8").Value = 137.452751
$ws.Range("O8").Value = 0.09035081822479055
$ws.Range("P8").Value = 0.09035081822479055
$ws.Range("Q8").Value = 30.40862628614634
$ws.Range("R8").Value = 273.677636575317
$ws.Range("S8").Value = 0.03028588158422349
$ws.Range("T8").Value = 0.03028588158422349
$ws.Range("A9").Value = "MuSCs"
$ws.Range("G9").Value = 0.663689
$ws.Range("H9").Value = 1.991067
$ws.Range("I9").Value = 0.3352031799963669
$ws.Range("J9").Value = 0.3352031799963669
$ws.Range("M9").Value = 270.6815896666666
$ws.Range("N9").Value = 812.044769
$ws.Range("O9").Value = 0.5337754885263156
$ws.Range("P9").Value = 0.5337754885263157
$ws.Range("Q9").Value = 179.6483935642803
$ws.Range("R9").Value = 1616.835542078523
$ws.Range("S9").Value = 0.1789232411581352
$ws.Range("T9").Value = 0.1789232411581353
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.246489
$ws.Range("H10").Value = 0.7394670000000001
$ws.Range("I10").Value = 0.1244918879687994
$ws.Range("J10").Value = 0.1244918879687994
$ws.Range("M10").Value = 52.63198466666667
$ws.Range("N10").Value = 157.895954
$ws.Range("O10").Value = 0.1037886003335349
$ws.Range("P10").Value = 0.1037886003335349
$ws.Range("Q10").Value = 12.973205268502
$ws.Range("R10").Value = 116.758847416518
$ws.Range("S10").Value = 0.01292083880516092
$ws.Range("T10").Value = 0.01292083880516092
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.246489
$ws.Range("H11").Value = 0.7394670000000001
$ws.Range("I11").Value = 0.1244918879687994
$ws.Range("J11").Value = 0.1244918879687994
$ws.Range("O11").Value = 0.2720850929153589
$ws.Range("P11").Value = 0.2720850929153589
$ws.Range("Q11").Value = 34.009667242327
$ws.Range("R11").Value = 306.087005180943
$ws.Range("S11").Value = 0.03387238690519923
$ws.Range("T11").Value = 0.03387238690519923
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.246489
$ws.Range("H12").Value = 0.7394670000000001
$ws.Range("I12").Value = 0.1244918879687994
$ws.Range("J12").Value = 0.1244918879687994
$ws.Range("M12").Value = 45.81758366666667
$ws.Range("N12").Value = 137.452751
$ws.Range("O12").Value = 0.09035081822479055
$ws.Range("P12").Value = 0.09035081822479055
$ws.Range("Q12").Value = 11.293530380413
$ws.Range("R12").Value = 101.641773423717
$ws.Range("S12").Value = 0.01124794394032998
$ws.Range("T12").Value = 0.01124794394032998
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.246489
$ws.Range("H13").Value = 0.7394670000000001
$ws.Range("I13").Value = 0.1244918879687994
$ws.Range("J13").Value = 0.1244918879687994
$ws.Range("M13").Value = 270.6815896666666
$ws.Range("N13").Value = 812.044769
$ws.Range("O13").Value = 0.5337754885263156
$ws.Range("P13").Value = 0.5337754885263157
$ws.Range("Q13").Value = 66.72003435534701
$ws.Range("R13").Value = 600.4803091981231
$ws.Range("S13").Value = 0.06645071831810925
$ws.Range("T13").Value = 0.06645071831810925
